$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet: rows 2-30 ---
$wsReco.Cells.Item(2, 1).Value = "BRVM - CONSOMMATION DE BASE     (**)"
$wsReco.Cells.Item(2, 2).Value = 0
$wsReco.Cells.Item(2, 3).Value = 3
$wsReco.Cells.Item(2, 4).Value = 785.39
$wsReco.Cells.Item(2, 5).Value = 267.46
$wsReco.Cells.Item(2, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(2, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(3, 1).Value = "BRVM-PRINCIPAL     (**)"
$wsReco.Cells.Item(3, 2).Value = 0
$wsReco.Cells.Item(3, 3).Value = 3
$wsReco.Cells.Item(3, 4).Value = 762.65
$wsReco.Cells.Item(3, 5).Value = 257.36
$wsReco.Cells.Item(3, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(3, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(4, 1).Value = "BRVM - INDUSTRIELS"
$wsReco.Cells.Item(4, 2).Value = 0
$wsReco.Cells.Item(4, 3).Value = 3
$wsReco.Cells.Item(4, 4).Value = 588.46
$wsReco.Cells.Item(4, 5).Value = 191.95
$wsReco.Cells.Item(4, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(4, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(5, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Cells.Item(5, 2).Value = 0
$wsReco.Cells.Item(5, 3).Value = 3
$wsReco.Cells.Item(5, 4).Value = 577.2
$wsReco.Cells.Item(5, 5).Value = 190.93
$wsReco.Cells.Item(5, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(5, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(6, 1).Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Cells.Item(6, 2).Value = 0
$wsReco.Cells.Item(6, 3).Value = 3
$wsReco.Cells.Item(6, 4).Value = 479.19
$wsReco.Cells.Item(6, 5).Value = 160.68
$wsReco.Cells.Item(6, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(6, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(7, 1).Value = "BRVM-PRESTIGE"
$wsReco.Cells.Item(7, 2).Value = 0
$wsReco.Cells.Item(7, 3).Value = 3
$wsReco.Cells.Item(7, 4).Value = 461.87
$wsReco.Cells.Item(7, 5).Value = 153.91
$wsReco.Cells.Item(7, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(7, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(8, 1).Value = "BRVM – COMPOSITE TOTAL RETURN     (**)"
$wsReco.Cells.Item(8, 2).Value = 0
$wsReco.Cells.Item(8, 3).Value = 3
$wsReco.Cells.Item(8, 4).Value = 440.62
$wsReco.Cells.Item(8, 5).Value = 147.6
$wsReco.Cells.Item(8, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(8, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(9, 1).Value = "BRVM - ENERGIE"
$wsReco.Cells.Item(9, 2).Value = 0
$wsReco.Cells.Item(9, 3).Value = 3
$wsReco.Cells.Item(9, 4).Value = 374.96
$wsReco.Cells.Item(9, 5).Value = 125.17
$wsReco.Cells.Item(9, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(9, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(10, 1).Value = "BRVM - SERVICES PUBLICS"
$wsReco.Cells.Item(10, 2).Value = 0
$wsReco.Cells.Item(10, 3).Value = 3
$wsReco.Cells.Item(10, 4).Value = 357.43
$wsReco.Cells.Item(10, 5).Value = 117.82
$wsReco.Cells.Item(10, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(10, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(11, 1).Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Cells.Item(11, 2).Value = 0
$wsReco.Cells.Item(11, 3).Value = 3
$wsReco.Cells.Item(11, 4).Value = 299.72
$wsReco.Cells.Item(11, 5).Value = 99.91
$wsReco.Cells.Item(11, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(11, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(12, 1).Value = "UNILEVER CI (UNLC)"
$wsReco.Cells.Item(12, 2).Value = 2
$wsReco.Cells.Item(12, 3).Value = 0
$wsReco.Cells.Item(12, 4).Value = 14.91
$wsReco.Cells.Item(12, 5).Value = 7.42
$wsReco.Cells.Item(12, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(12, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(13, 1).Value = "SETAO CI (STAC)"
$wsReco.Cells.Item(13, 2).Value = 2
$wsReco.Cells.Item(13, 3).Value = 1
$wsReco.Cells.Item(13, 4).Value = 7.73
$wsReco.Cells.Item(13, 5).Value = -6.83
$wsReco.Cells.Item(13, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(13, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(14, 1).Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Cells.Item(14, 2).Value = 1
$wsReco.Cells.Item(14, 3).Value = 0
$wsReco.Cells.Item(14, 4).Value = 7.45
$wsReco.Cells.Item(14, 5).Value = 7.45
$wsReco.Cells.Item(14, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(14, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(15, 1).Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$wsReco.Cells.Item(15, 2).Value = 2
$wsReco.Cells.Item(15, 3).Value = 1
$wsReco.Cells.Item(15, 4).Value = 7.37
$wsReco.Cells.Item(15, 5).Value = -7.37
$wsReco.Cells.Item(15, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(15, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(16, 1).Value = "UNIWAX CI (UNXC)"
$wsReco.Cells.Item(16, 2).Value = 2
$wsReco.Cells.Item(16, 3).Value = 1
$wsReco.Cells.Item(16, 4).Value = 7.15
$wsReco.Cells.Item(16, 5).Value = -7.46
$wsReco.Cells.Item(16, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(16, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(17, 1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$wsReco.Cells.Item(17, 2).Value = 2
$wsReco.Cells.Item(17, 3).Value = 1
$wsReco.Cells.Item(17, 4).Value = 7.03
$wsReco.Cells.Item(17, 5).Value = -7.37
$wsReco.Cells.Item(17, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(17, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(18, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$wsReco.Cells.Item(18, 2).Value = 1
$wsReco.Cells.Item(18, 3).Value = 0
$wsReco.Cells.Item(18, 4).Value = 3.85
$wsReco.Cells.Item(18, 5).Value = 3.85
$wsReco.Cells.Item(18, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(18, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(19, 1).Value = "SOLIBRA CI (SLBC)"
$wsReco.Cells.Item(19, 2).Value = 1
$wsReco.Cells.Item(19, 3).Value = 0
$wsReco.Cells.Item(19, 4).Value = 2.17
$wsReco.Cells.Item(19, 5).Value = 2.17
$wsReco.Cells.Item(19, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(19, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(20, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$wsReco.Cells.Item(20, 2).Value = 1
$wsReco.Cells.Item(20, 3).Value = 0
$wsReco.Cells.Item(20, 4).Value = 1.65
$wsReco.Cells.Item(20, 5).Value = 1.65
$wsReco.Cells.Item(20, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(20, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(21, 1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$wsReco.Cells.Item(21, 2).Value = 0
$wsReco.Cells.Item(21, 3).Value = 1
$wsReco.Cells.Item(21, 4).Value = -1.59
$wsReco.Cells.Item(21, 5).Value = -1.59
$wsReco.Cells.Item(21, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(21, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(22, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$wsReco.Cells.Item(22, 2).Value = 0
$wsReco.Cells.Item(22, 3).Value = 1
$wsReco.Cells.Item(22, 4).Value = -1.64
$wsReco.Cells.Item(22, 5).Value = -1.64
$wsReco.Cells.Item(22, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(22, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(23, 1).Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$wsReco.Cells.Item(23, 2).Value = 0
$wsReco.Cells.Item(23, 3).Value = 1
$wsReco.Cells.Item(23, 4).Value = -2
$wsReco.Cells.Item(23, 5).Value = -2
$wsReco.Cells.Item(23, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(23, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(24, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Cells.Item(24, 2).Value = 0
$wsReco.Cells.Item(24, 3).Value = 1
$wsReco.Cells.Item(24, 4).Value = -2.29
$wsReco.Cells.Item(24, 5).Value = -2.29
$wsReco.Cells.Item(24, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(24, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(25, 1).Value = "BANK OF AFRICA BF (BOABF)"
$wsReco.Cells.Item(25, 2).Value = 0
$wsReco.Cells.Item(25, 3).Value = 1
$wsReco.Cells.Item(25, 4).Value = -2.72
$wsReco.Cells.Item(25, 5).Value = -2.72
$wsReco.Cells.Item(25, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(25, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(26, 1).Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Cells.Item(26, 2).Value = 0
$wsReco.Cells.Item(26, 3).Value = 1
$wsReco.Cells.Item(26, 4).Value = -2.78
$wsReco.Cells.Item(26, 5).Value = -2.78
$wsReco.Cells.Item(26, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(26, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(27, 1).Value = "BERNABE CI (BNBC)"
$wsReco.Cells.Item(27, 2).Value = 0
$wsReco.Cells.Item(27, 3).Value = 1
$wsReco.Cells.Item(27, 4).Value = -2.94
$wsReco.Cells.Item(27, 5).Value = -2.94
$wsReco.Cells.Item(27, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(27, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(28, 1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$wsReco.Cells.Item(28, 2).Value = 1
$wsReco.Cells.Item(28, 3).Value = 2
$wsReco.Cells.Item(28, 4).Value = -3.98
$wsReco.Cells.Item(28, 5).Value = 1.79
$wsReco.Cells.Item(28, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(28, 7).Value = "👀 À surveiller"

$wsReco.Cells.Item(29, 1).Value = "FILTISAC CI (FTSC)"
$wsReco.Cells.Item(29, 2).Value = 0
$wsReco.Cells.Item(29, 3).Value = 1
$wsReco.Cells.Item(29, 4).Value = -7.26
$wsReco.Cells.Item(29, 5).Value = -7.26
$wsReco.Cells.Item(29, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(29, 7).Value = "➖ Neutre"

$wsReco.Cells.Item(30, 1).Value = "SICABLE CI (CABC)"
$wsReco.Cells.Item(30, 2).Value = 0
$wsReco.Cells.Item(30, 3).Value = 1
$wsReco.Cells.Item(30, 4).Value = -7.5
$wsReco.Cells.Item(30, 5).Value = -7.5
$wsReco.Cells.Item(30, 6).Value = "🟡 Observer"
$wsReco.Cells.Item(30, 7).Value = "➖ Neutre"

# --- Top_YTD sheet: rows 2-11 (column B only changes) ---
$wsYtd.Cells.Item(2, 2).Value = 4634.67
$wsYtd.Cells.Item(3, 2).Value = 4343.93
$wsYtd.Cells.Item(4, 2).Value = 2497.03
$wsYtd.Cells.Item(5, 2).Value = 2399.89
$wsYtd.Cells.Item(6, 2).Value = 1652.11
$wsYtd.Cells.Item(7, 2).Value = 1537.87
$wsYtd.Cells.Item(8, 2).Value = 1404.59
$wsYtd.Cells.Item(9, 2).Value = 1038.84
$wsYtd.Cells.Item(10, 2).Value = 952.38
$wsYtd.Cells.Item(11, 2).Value = 698.88
